$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FoTOMRAEL")

# AEO updates - Endogenous learning for batteries:
# "Share of Technology Outside Modeled Region" revised from 25% to 90%.
$ws.Range("B2").Value = 0.9

# Make the FoTOMRAEL sheet the active/selected tab with B3 selected,
# matching the saved view state in the updated workbook.
$ws.Activate()
$ws.Range("B3").Select() | Out-Null
